$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 9
$ws.Range("B5").Value = "testuser@gmail.com"
$ws.Range("C5").Value = "ss"
$ws.Range("D5").Value = "ss"
$ws.Range("E5").Value = "ss"
$ws.Range("F5").Value = "Customer"

$ws.Columns.AutoFit() | Out-Null
